$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 18.14140940046379
$ws.Range("C2").Value = 11.74475127611383
$ws.Range("D2").Value = 9.914670089302165
$ws.Range("F2").Value = 30.35817781966811
$ws.Range("G2").Value = 29.89473087637171
$ws.Range("H2").Value = 14.67751478693199
$ws.Range("J2").Value = 10.46819755913964
$ws.Range("L2").Value = 11.66141587602543
$ws.Range("O2").Value = 22.47356275402743
$ws.Range("B3").Value = 17.49512930106765
$ws.Range("C3").Value = 11.56519436279091
$ws.Range("D3").Value = 9.914261041400732
$ws.Range("F3").Value = 30.4943709402507
$ws.Range("G3").Value = 30.06292800558843
$ws.Range("H3").Value = 14.7442300333797
$ws.Range("J3").Value = 10.50328319840692
$ws.Range("L3").Value = 11.6245366695048
$ws.Range("O3").Value = 22.59231876496447
$ws.Range("B4").Value = 17.08616403992307
$ws.Range("C4").Value = 11.45413522548096
$ws.Range("D4").Value = 9.915284748456017
$ws.Range("F4").Value = 30.58628715720752
$ws.Range("G4").Value = 30.17779478386595
$ws.Range("H4").Value = 14.78797393332698
$ws.Range("J4").Value = 10.52595140273449
$ws.Range("L4").Value = 11.60302834398714
$ws.Range("O4").Value = 22.67098359697485
$ws.Range("B5").Value = 16.91667546289918
$ws.Range("C5").Value = 11.40871999684409
$ws.Range("D5").Value = 9.916023546687883
$ws.Range("F5").Value = 30.62582267168682
$ws.Range("G5").Value = 30.22750111956405
$ws.Range("H5").Value = 14.80649909307044
$ws.Range("J5").Value = 10.53547269181328
$ws.Range("L5").Value = 11.59455425153929
$ws.Range("O5").Value = 22.70448246725735
$ws.Range("B6").Value = 16.88836792472812
$ws.Range("C6").Value = 11.40117066782581
$ws.Range("D6").Value = 9.916165678151714
$ws.Range("F6").Value = 30.63251288359834
$ws.Range("G6").Value = 30.23592926964983
$ws.Range("H6").Value = 14.80961740988206
$ws.Range("J6").Value = 10.53707085970758
$ws.Range("L6").Value = 11.59316484002602
$ws.Range("O6").Value = 22.71013195096581
$ws.Range("B7").Value = 17.08388941323213
$ws.Range("C7").Value = 11.45352331702556
$ws.Range("D7").Value = 9.915293408552344
$ws.Range("F7").Value = 30.58681193819371
$ws.Range("G7").Value = 30.17845343431514
$ws.Range("H7").Value = 14.78822093921447
$ws.Range("J7").Value = 10.52607865996241
$ws.Range("L7").Value = 11.60291287564064
$ws.Range("O7").Value = 22.67142953868719
$ws.Range("B8").Value = 17.92121949698125
$ws.Range("C8").Value = 11.6830356408769
$ws.Range("D8").Value = 9.914265109866131
$ws.Range("F8").Value = 30.40341213725929
$ws.Range("G8").Value = 29.95030840661072
$ws.Range("H8").Value = 14.69994111181612
$ws.Range("J8").Value = 10.4800619774831
$ws.Range("L8").Value = 11.64846737826615
$ws.Range("O8").Value = 22.51331492326075
$ws.Range("B9").Value = 19.45820809253776
$ws.Range("C9").Value = 12.12470245141216
$ws.Range("D9").Value = 9.922314945266463
$ws.Range("F9").Value = 30.10985049501152
$ws.Range("G9").Value = 29.59568829588758
$ws.Range("H9").Value = 14.54888657293744
$ws.Range("J9").Value = 10.39871682171028
$ws.Range("L9").Value = 11.74657572878298
$ws.Range("O9").Value = 22.24900101952358
$ws.Range("B10").Value = 20.51335592316653
$ws.Range("C10").Value = 12.4415094693615
$ws.Range("D10").Value = 9.934295706465814
$ws.Range("F10").Value = 29.9348357826732
$ws.Range("G10").Value = 29.39275057571173
$ws.Range("H10").Value = 14.45135646052431
$ws.Range("J10").Value = 10.34432243261253
$ws.Range("L10").Value = 11.82367979163881
$ws.Range("O10").Value = 22.0828901633131
$ws.Range("B11").Value = 20.97555491419136
$ws.Range("C11").Value = 12.58345037846499
$ws.Range("D11").Value = 9.941046432926017
$ws.Range("F11").Value = 29.86412478439616
$ws.Range("G11").Value = 29.31315015653858
$ws.Range("H11").Value = 14.40990852843186
$ws.Range("J11").Value = 10.3207320230515
$ws.Range("L11").Value = 11.85977361913809
$ws.Range("O11").Value = 22.01346116715323
$ws.Range("B12").Value = 21.14789841991029
$ws.Range("C12").Value = 12.6368460673621
$ws.Range("D12").Value = 9.943788157111646
$ws.Range("F12").Value = 29.83863492503324
$ws.Range("G12").Value = 29.28485294266675
$ws.Range("H12").Value = 14.39463320923325
$ws.Range("J12").Value = 10.31196404984409
$ws.Range("L12").Value = 11.87358134044042
$ws.Range("O12").Value = 21.98805611523642
$ws.Range("B13").Value = 21.11090219418308
$ws.Range("C13").Value = 12.6253627436381
$ws.Range("D13").Value = 9.943189460831187
$ws.Range("F13").Value = 29.84406728821661
$ws.Range("G13").Value = 29.29086490015247
$ws.Range("H13").Value = 14.39790433180278
$ws.Range("J13").Value = 10.31384505374139
$ws.Range("L13").Value = 11.8706014931698
$ws.Range("O13").Value = 21.99348807128205
$ws.Range("B14").Value = 20.98978809960123
$ws.Range("C14").Value = 12.58785060363051
$ws.Range("D14").Value = 9.94126829097976
$ws.Range("F14").Value = 29.86200189145836
$ws.Range("G14").Value = 29.31078504117542
$ws.Range("H14").Value = 14.40864339685411
$ws.Range("J14").Value = 10.32000736982025
$ws.Range("L14").Value = 11.86090682659151
$ws.Range("O14").Value = 22.01135329500545
$ws.Range("B15").Value = 20.91524967882819
$ws.Range("C15").Value = 12.56482600646688
$ws.Range("D15").Value = 9.940115608757537
$ws.Range("F15").Value = 29.87315512925376
$ws.Range("G15").Value = 29.32322756658514
$ws.Range("H15").Value = 14.41527610667861
$ws.Range("J15").Value = 10.32380345932324
$ws.Range("L15").Value = 11.85498656318723
$ws.Range("O15").Value = 22.02241178623174
$ws.Range("B16").Value = 20.48278124189433
$ws.Range("C16").Value = 12.43218593176972
$ws.Range("D16").Value = 9.933880563130494
$ws.Range("F16").Value = 29.93963664312136
$ws.Range("G16").Value = 29.39821002142552
$ws.Range("H16").Value = 14.45412394867826
$ws.Range("J16").Value = 10.34588727829284
$ws.Range("L16").Value = 11.82134090466173
$ws.Range("O16").Value = 22.08755130968184
$ws.Range("B17").Value = 20.21282623517714
$ws.Range("C17").Value = 12.35022817813506
$ws.Range("D17").Value = 9.930387538363513
$ws.Range("F17").Value = 29.98270611419626
$ws.Range("G17").Value = 29.44747879551798
$ws.Range("H17").Value = 14.47870371884882
$ws.Range("J17").Value = 10.35973000391522
$ws.Range("L17").Value = 11.80095658966508
$ws.Range("O17").Value = 22.12908655638329
$ws.Range("B18").Value = 20.05588673890565
$ws.Range("C18").Value = 12.30288576852642
$ws.Range("D18").Value = 9.928500968986995
$ws.Range("F18").Value = 30.00831625467118
$ws.Range("G18").Value = 29.47701289328056
$ws.Range("H18").Value = 14.4931160796135
$ws.Range("J18").Value = 10.36780062411152
$ws.Range("L18").Value = 11.78932837956123
$ws.Range("O18").Value = 22.15355381464307
$ws.Range("B19").Value = 20.00246715395293
$ws.Range("C19").Value = 12.28682295595189
$ws.Range("D19").Value = 9.927883301754511
$ws.Range("F19").Value = 30.01713108685006
$ws.Range("G19").Value = 29.48721751522045
$ws.Range("H19").Value = 14.4980430313503
$ws.Range("J19").Value = 10.37055188039443
$ws.Range("L19").Value = 11.7854080107528
$ws.Range("O19").Value = 22.16193703834668
$ws.Range("B20").Value = 20.24173703066387
$ws.Range("C20").Value = 12.35897398183142
$ws.Range("D20").Value = 9.930746707500926
$ws.Range("F20").Value = 29.97803454706695
$ws.Range("G20").Value = 29.44211014297847
$ws.Range("H20").Value = 14.47605872457396
$ws.Range("J20").Value = 10.35824518246824
$ws.Range("L20").Value = 11.80311661541161
$ws.Range("O20").Value = 22.12460528365174
$ws.Range("B21").Value = 21.02543591094002
$ws.Range("C21").Value = 12.59887876269585
$ws.Range("D21").Value = 9.941827567948424
$ws.Range("F21").Value = 29.85669908826476
$ws.Range("G21").Value = 29.30488378762517
$ws.Range("H21").Value = 14.40547766834395
$ws.Range("J21").Value = 10.3181928703753
$ws.Range("L21").Value = 11.86375064627112
$ws.Range("O21").Value = 22.00608176326341
$ws.Range("B22").Value = 21.52195425473571
$ws.Range("C22").Value = 12.75358773485035
$ws.Range("D22").Value = 9.950149190930734
$ws.Range("F22").Value = 29.7849030037151
$ws.Range("G22").Value = 29.22596421139951
$ws.Range("H22").Value = 14.36179779715629
$ws.Range("J22").Value = 10.29297896758686
$ws.Range("L22").Value = 11.9041899513632
$ws.Range("O22").Value = 21.93378708091132
$ws.Range("B23").Value = 21.2584237253133
$ws.Range("C23").Value = 12.67122036713452
$ws.Range("D23").Value = 9.945609571045258
$ws.Range("F23").Value = 29.82253326199502
$ws.Range("G23").Value = 29.26709455254387
$ws.Range("H23").Value = 14.38488636198601
$ws.Range("J23").Value = 10.30634826114856
$ws.Range("L23").Value = 11.88253475101733
$ws.Range("O23").Value = 21.97189801145809
$ws.Range("B24").Value = 20.22867186643314
$ws.Range("C24").Value = 12.35502069793155
$ws.Range("D24").Value = 9.930583948196352
$ws.Range("F24").Value = 29.98014391869435
$ws.Range("G24").Value = 29.44453354667599
$ws.Range("H24").Value = 14.47725365081761
$ws.Range("J24").Value = 10.35891612052683
$ws.Range("L24").Value = 11.80213978380207
$ws.Range("O24").Value = 22.12662943546187
$ws.Range("B25").Value = 19.05477197518198
$ws.Range("C25").Value = 12.00638957754411
$ws.Range("D25").Value = 9.919066419287013
$ws.Range("F25").Value = 30.18215219574993
$ws.Range("G25").Value = 29.68158074521383
$ws.Range("H25").Value = 14.58738903843093
$ws.Range("J25").Value = 10.41977611538033
$ws.Range("L25").Value = 11.74657572878298
$ws.Range("O25").Value = 22.31558764326136
